$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: translate to snake_case column names
$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'

# Fix A190: drop stray _x000D_ / CR artifact and title-case "de" -> "De"
$ws.Range("A190").Value = 'Estado De México'

# Title-case the Spanish linking words (de/del/el/la/las/los/y) in municipality/state names
$ws.Range("B8").Value = 'Pabellón De Arteaga'
$ws.Range("B9").Value = 'Rincón De Romos'
$ws.Range("B10").Value = 'San Francisco De Los Romo'
$ws.Range("B15").Value = 'Playas De Rosarito'
$ws.Range("B35").Value = 'Chiapa De Corzo'
$ws.Range("B57").Value = 'San Cristóbal De Las Casas'
$ws.Range("B92").Value = 'Guadalupe Y Calvo'
$ws.Range("B94").Value = 'Hidalgo Del Parral'
$ws.Range("B110").Value = 'San Francisco De Borja'
$ws.Range("B111").Value = 'San Francisco Del Oro'
$ws.Range("B117").Value = 'Valle De Zaragoza'
$ws.Range("B139").Value = 'Villa De Álvarez'
$ws.Range("A141").Value = 'Ciudad De México'
$ws.Range("B145").Value = 'Cuajimalpa De Morelos'
$ws.Range("B160").Value = 'Coneto De Comonfort'
$ws.Range("B170").Value = 'Nombre De Dios'
$ws.Range("B173").Value = 'Pánuco De Coronado'
$ws.Range("B179").Value = 'San Juan Del Río'
$ws.Range("B180").Value = 'San Luis Del Cordero'
$ws.Range("B181").Value = 'San Pedro Del Gallo'
$ws.Range("B190").Value = 'Acambay De Ruíz Castañeda'
$ws.Range("B192").Value = 'Almoloya De Alquisiras'
$ws.Range("B193").Value = 'Almoloya De Juárez'
$ws.Range("B199").Value = 'Atizapán De Zaragoza'
$ws.Range("B205").Value = 'Chapa De Mota'
$ws.Range("B208").Value = 'Coacalco De Berriozábal'
$ws.Range("B213").Value = 'Ecatepec De Morelos'
$ws.Range("B218").Value = 'Ixtapan De La Sal'
$ws.Range("B231").Value = 'Naucalpan De Juárez'
$ws.Range("B241").Value = 'San Felipe Del Progreso'
$ws.Range("B242").Value = 'San Martín De Las Pirámides'
$ws.Range("B244").Value = 'San Simón De Guerrero'
$ws.Range("B245").Value = 'Soyaniquilpan De Juárez'
$ws.Range("B254").Value = 'Tenango Del Valle'
$ws.Range("B265").Value = 'Tlalnepantla De Baz'
$ws.Range("B271").Value = 'Valle De Bravo'
$ws.Range("B272").Value = 'Valle De Chalco Solidaridad'
$ws.Range("B273").Value = 'Villa De Allende'
$ws.Range("B283").Value = 'San Miguel De Allende'
$ws.Range("B284").Value = 'Apaseo El Alto'
$ws.Range("B285").Value = 'Apaseo El Grande'
$ws.Range("B291").Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range("B295").Value = 'Jaral Del Progreso'
$ws.Range("B303").Value = 'Purísima Del Rincón'
$ws.Range("B307").Value = 'San Diego De La Unión'
$ws.Range("B309").Value = 'San Francisco Del Rincón'
$ws.Range("B311").Value = 'San Luis De La Paz'
$ws.Range("B313").Value = 'Santa Cruz De Juventino Rosas'
$ws.Range("B315").Value = 'Silao De La Victoria'
$ws.Range("B319").Value = 'Valle De Santiago'
$ws.Range("B324").Value = 'Acapulco De Juárez'
$ws.Range("B327").Value = 'Ajuchitlán Del Progreso'
$ws.Range("B328").Value = 'Alcozauca De Guerrero'
$ws.Range("B332").Value = 'Atenango Del Río'
$ws.Range("B333").Value = 'Atoyac De Álvarez'
$ws.Range("B334").Value = 'Ayutla De Los Libres'
$ws.Range("B337").Value = 'Buenavista De Cuéllar'
$ws.Range("B338").Value = 'Chilapa De Álvarez'
$ws.Range("B339").Value = 'Chilpancingo De Los Bravo'
$ws.Range("B340").Value = 'Coahuayutla De José María Izazaga'
$ws.Range("B345").Value = 'Coyuca De Benítez'
$ws.Range("B346").Value = 'Coyuca De Catalán'
$ws.Range("B349").Value = 'Cuetzala Del Progreso'
$ws.Range("B350").Value = 'Cutzamala De Pinzón'
$ws.Range("B355").Value = 'Huitzuco De Los Figueroa'
$ws.Range("B356").Value = 'Iguala De La Independencia'
$ws.Range("B358").Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range("B361").Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range("B364").Value = 'Mártir De Cuilapan'
$ws.Range("B377").Value = 'Taxco De Alarcón'
$ws.Range("B379").Value = 'Técpan De Galeana'
$ws.Range("B381").Value = 'Tepecoacuilco De Trujano'
$ws.Range("B383").Value = 'Tixtla De Guerrero'
$ws.Range("B386").Value = 'Tlalixtaquilla De Maldonado'
$ws.Range("B387").Value = 'Tlapa De Comonfort'
$ws.Range("B398").Value = 'Agua Blanca De Iturbide'
$ws.Range("B404").Value = 'Atotonilco De Tula'
$ws.Range("B405").Value = 'Atotonilco El Grande'
$ws.Range("B410").Value = 'Cuautepec De Hinojosa'
$ws.Range("B414").Value = 'Huasca De Ocampo'
$ws.Range("B417").Value = 'Huejutla De Reyes'
$ws.Range("B422").Value = 'Mineral Del Monte'
$ws.Range("B423").Value = 'Mixquiahuala De Juárez'
$ws.Range("B424").Value = 'Pachuca De Soto'
$ws.Range("B426").Value = 'Progreso De Obregón'
$ws.Range("B429").Value = 'Santiago De Anaya'
$ws.Range("B430").Value = 'Santiago Tulantepec De Lugo Guerrero'
$ws.Range("B435").Value = 'Tepehuacán De Guerrero'
$ws.Range("B436").Value = 'Tepeji Del Río De Ocampo'
$ws.Range("B438").Value = 'Tezontepec De Aldama'
$ws.Range("B443").Value = 'Tula De Allende'
$ws.Range("B444").Value = 'Tulancingo De Bravo'
$ws.Range("B445").Value = 'Villa De Tezontepec'
$ws.Range("B447").Value = 'Zacualtipán De Ángeles'
$ws.Range("B448").Value = 'Zapotlán De Juárez'
$ws.Range("B453").Value = 'Acatlán De Juárez'
$ws.Range("B454").Value = 'Ahualulco De Mercado'
$ws.Range("B458").Value = 'Atemajac De Brizuela'
$ws.Range("B461").Value = 'Atotonilco El Alto'
$ws.Range("B463").Value = 'Autlán De Navarro'
$ws.Range("B469").Value = 'Cañadas De Obregón'
$ws.Range("B476").Value = 'Cuautitlán De García Barragán'
$ws.Range("B481").Value = 'Encarnación De Díaz'
$ws.Range("B488").Value = 'Huejuquilla El Alto'
$ws.Range("B489").Value = 'Ixtlahuacán De Los Membrillos'
$ws.Range("B490").Value = 'Ixtlahuacán Del Río'
$ws.Range("B494").Value = 'Jilotlán De Los Dolores'
$ws.Range("B500").Value = 'La Manzanilla De La Paz'
$ws.Range("B501").Value = 'Lagos De Moreno'
$ws.Range("B508").Value = 'Ojuelos De Jalisco'
$ws.Range("B513").Value = 'San Cristóbal De La Barranca'
$ws.Range("B514").Value = 'San Diego De Alejandría'
$ws.Range("B516").Value = 'San Juan De Los Lagos'
$ws.Range("B517").Value = 'San Juanito De Escobedo'
$ws.Range("B520").Value = 'San Martín De Bolaños'
$ws.Range("B522").Value = 'San Miguel El Alto'
$ws.Range("B523").Value = 'San Sebastián Del Oeste'
$ws.Range("B524").Value = 'Santa María De Los Ángeles'
$ws.Range("B527").Value = 'Talpa De Allende'
$ws.Range("B528").Value = 'Tamazula De Gordiano'
$ws.Range("B534").Value = 'Teocuitatlán De Corona'
$ws.Range("B535").Value = 'Tepatitlán De Morelos'
$ws.Range("B538").Value = 'Tizapán El Alto'
$ws.Range("B539").Value = 'Tlajomulco De Zúñiga'
$ws.Range("B550").Value = 'Unión De San Antonio'
$ws.Range("B551").Value = 'Unión De Tula'
$ws.Range("B552").Value = 'Valle De Guadalupe'
$ws.Range("B553").Value = 'Valle De Juárez'
$ws.Range("B558").Value = 'Zacoalco De Torres'
$ws.Range("B561").Value = 'Zapotitlán De Vadillo'
$ws.Range("B562").Value = 'Zapotlán Del Rey'
$ws.Range("B563").Value = 'Zapotlán El Grande'
$ws.Range("B586").Value = 'Coalcomán De Vázquez Pallares'
$ws.Range("B588").Value = 'Cojumatlán De Régules'
$ws.Range("B652").Value = 'Tiquicheo De Nicolás Romero'
$ws.Range("B676").Value = 'Coatlán Del Río'
$ws.Range("B687").Value = 'Puente De Ixtla'
$ws.Range("B693").Value = 'Tlaltizapán De Zapata'
$ws.Range("B700").Value = 'Zacualpan De Amilpas'
$ws.Range("B704").Value = 'Amatlán De Cañas'
$ws.Range("B705").Value = 'Bahía De Banderas'
$ws.Range("B709").Value = 'Ixtlán Del Río'
$ws.Range("B716").Value = 'Santa María Del Oro'
$ws.Range("B731").Value = 'San Nicolás De Los Garza'
$ws.Range("B734").Value = 'Acatlán De Pérez Figueroa'
$ws.Range("B736").Value = 'Ciénega De Zimatlán'
$ws.Range("B738").Value = 'Constancia Del Rosario'
$ws.Range("B740").Value = 'Cuyamecalco Villa De Zaragoza'
$ws.Range("B741").Value = 'Guadalupe De Ramírez'
$ws.Range("B742").Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range("B743").Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range("B744").Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range("B745").Value = 'Huautla De Jiménez'
$ws.Range("B746").Value = 'Ixtlán De Juárez'
$ws.Range("B747").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B753").Value = 'Mariscala De Juárez'
$ws.Range("B754").Value = 'Mártires De Tacubaya'
$ws.Range("B757").Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range("B759").Value = 'Oaxaca De Juárez'
$ws.Range("B760").Value = 'Ocotlán De Morelos'
$ws.Range("B761").Value = 'Putla Villa De Guerrero'
$ws.Range("B785").Value = 'San Juan Bautista Lo De Soto'
$ws.Range("B813").Value = 'San Pedro Y San Pablo Teposcolula'
$ws.Range("B827").Value = 'Santa Lucía Del Camino'
$ws.Range("B833").Value = 'Santa María Jalapa Del Marqués'
$ws.Range("B854").Value = 'Santo Domingo De Morelos'
$ws.Range("B859").Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range("B860").Value = 'Teotitlán De Flores Magón'
$ws.Range("B862").Value = 'Tlacolula De Matamoros'
$ws.Range("B864").Value = 'Villa De Chilapa De Díaz'
$ws.Range("B865").Value = 'Villa De Etla'
$ws.Range("B866").Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range("B867").Value = 'Villa Sola De Vega'
$ws.Range("B870").Value = 'Zimatlán De Álvarez'
$ws.Range("B883").Value = 'Ayotoxco De Guerrero'
$ws.Range("B886").Value = 'Chalchicomula De Sesma'
$ws.Range("B894").Value = 'Chila De La Sal'
$ws.Range("B901").Value = 'Cuayuca De Andrade'
$ws.Range("B910").Value = 'Huehuetlán El Chico'
$ws.Range("B913").Value = 'Izúcar De Matamoros'
$ws.Range("B920").Value = 'Los Reyes De Juárez'
$ws.Range("B926").Value = 'Palmar De Bravo'
$ws.Range("B941").Value = 'San Nicolás De Los Ranchos'
$ws.Range("B943").Value = 'San Salvador El Seco'
$ws.Range("B948").Value = 'Tecali De Herrera'
$ws.Range("B954").Value = 'Tepanco De López'
$ws.Range("B955").Value = 'Tepatlaxco De Hidalgo'
$ws.Range("B959").Value = 'Tepexi De Rodríguez'
$ws.Range("B961").Value = 'Tetela De Ocampo'
$ws.Range("B966").Value = 'Tlacotepec De Benito Juárez'
$ws.Range("B978").Value = 'Xayacatlán De Bravo'
$ws.Range("B989").Value = 'Amealco De Bonfil'
$ws.Range("B990").Value = 'Cadereyta De Montes'
$ws.Range("B993").Value = 'Landa De Matamoros'
$ws.Range("B995").Value = 'Pinal De Amoles'
$ws.Range("B997").Value = 'San Juan Del Río'
$ws.Range("B1007").Value = 'Armadillo De Los Infante'
$ws.Range("B1011").Value = 'Ciudad Del Maíz'
$ws.Range("B1017").Value = 'Mexquitic De Carmona'
$ws.Range("B1022").Value = 'San Ciro De Acosta'
$ws.Range("B1025").Value = 'Soledad De Graciano Sánchez'
$ws.Range("B1028").Value = 'Tanquián De Escobedo'
$ws.Range("B1031").Value = 'Villa De Arriaga'
$ws.Range("B1032").Value = 'Villa De Ramos'
$ws.Range("B1078").Value = 'Nacozari De García'
$ws.Range("B1096").Value = 'Jalpa De Méndez'
$ws.Range("B1110").Value = 'Soto La Marina'
$ws.Range("B1120").Value = 'Contla De Juan Cuamatzi'
$ws.Range("B1124").Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range("B1125").Value = 'Nanacamilpa De Mariano Arista'
$ws.Range("B1128").Value = 'Sanctórum De Lázaro Cárdenas'
$ws.Range("B1145").Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range("B1148").Value = 'Amatlán De Los Reyes'
$ws.Range("B1153").Value = 'Boca Del Río'
$ws.Range("B1154").Value = 'Camarón De Tejeda'
$ws.Range("B1156").Value = 'Cazones De Herrera'
$ws.Range("B1164").Value = 'Cosamaloapan De Carpio'
$ws.Range("B1177").Value = 'Hueyapan De Ocampo'
$ws.Range("B1178").Value = 'Huiloapan De Cuauhtémoc'
$ws.Range("B1179").Value = 'Ignacio De La Llave'
$ws.Range("B1181").Value = 'Ixhuatlán De Madero'
$ws.Range("B1189").Value = 'Juchique De Ferrer'
$ws.Range("B1192").Value = 'Lerdo De Tejada'
$ws.Range("B1195").Value = 'Martínez De La Torre'
$ws.Range("B1197").Value = 'Medellín De Bravo'
$ws.Range("B1200").Value = 'Nanchital De Lázaro Cárdenas Del Río'
$ws.Range("B1208").Value = 'Paso Del Macho'
$ws.Range("B1210").Value = 'Poza Rica De Hidalgo'
$ws.Range("B1239").Value = 'Vega De Alatorre'
$ws.Range("B1265").Value = 'Cañitas De Felipe Pescador'
$ws.Range("B1267").Value = 'Concepción Del Oro'
$ws.Range("B1282").Value = 'Mezquital Del Oro'
$ws.Range("B1286").Value = 'Nochistlán De Mejía'
$ws.Range("B1287").Value = 'Noria De Ángeles'
$ws.Range("B1297").Value = 'Teúl De González Ortega'
$ws.Range("B1298").Value = 'Tlaltenango De Sánchez Román'
$ws.Range("B1300").Value = 'Villa De Cos'

# Tiny 1-ULP float recalculation fix for rows where C = 12 (12/12759)
$ws.Range("D91").Value = 0.0009405125793557488
$ws.Range("D209").Value = 0.0009405125793557488
$ws.Range("D246").Value = 0.0009405125793557488
$ws.Range("D298").Value = 0.0009405125793557488
$ws.Range("D405").Value = 0.0009405125793557488
$ws.Range("D649").Value = 0.0009405125793557488
$ws.Range("D683").Value = 0.0009405125793557488
$ws.Range("D709").Value = 0.0009405125793557488
$ws.Range("D862").Value = 0.0009405125793557488
$ws.Range("D874").Value = 0.0009405125793557488
$ws.Range("D913").Value = 0.0009405125793557488
$ws.Range("D1004").Value = 0.0009405125793557488
$ws.Range("D1037").Value = 0.0009405125793557488
$ws.Range("D1053").Value = 0.0009405125793557488
$ws.Range("D1161").Value = 0.0009405125793557488

# Remove trailing footer/metadata rows (1308:1313) so the used range ends at row 1307
$ws.Range("A1308:A1313").EntireRow.Delete()

